# "Commit for entire setup and testmodule change"
#
# The valid-login test row's e-mail credential (sheet "testCaseData",
# cell B2) is swapped for a new throwaway address. Excel rebuilds the
# shared-string table on save, so every other <si> naturally renumbers -
# setting the cell Value is enough to reproduce that ripple.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "gowtham@yahoo.com"

# The saved worksheet view's selection ends up anchored on B16 (the user
# ctrl-selected B2 and B16, leaving B16 as the active cell).
$ws.Range("B16").Select()
